# Applies crypto price/volume/coin updates per commit "Updated symbol list on Mon Jan  2 23:36:10 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a literal text value, preserving text type even
# when the value looks numeric (e.g. "246.09", "0.86%"), without leaving
# behind a stray explicit cell style.
function Set-TextValue($range, $val) {
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '246.09'
Set-TextValue $ws.Range('E2') '0.86%'
Set-TextValue $ws.Range('D3') '29.34'
Set-TextValue $ws.Range('E3') '6.63%'
Set-TextValue $ws.Range('D4') '5.180'
Set-TextValue $ws.Range('E4') '0.93%'
Set-TextValue $ws.Range('D5') '0.05737'
Set-TextValue $ws.Range('E5') '0.73%'
Set-TextValue $ws.Range('D6') '6.568'
Set-TextValue $ws.Range('E6') '0.71%'
$ws.Range('B7').Value = 'GateToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws.Range('D7') '3.109'
Set-TextValue $ws.Range('E7') '3.30%'
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range('D8') '0.8591'
Set-TextValue $ws.Range('E8') '4.78%'
$ws.Range('B9').Value = 'FTXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue $ws.Range('D9') '0.8625'
Set-TextValue $ws.Range('E9') '-0.68%'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws.Range('D10') '0.1363'
Set-TextValue $ws.Range('E10') '2.32%'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws.Range('D11') '0.07087'
Set-TextValue $ws.Range('E11') '1.92%'
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws.Range('D12') '0.03058'
Set-TextValue $ws.Range('E12') '6.39%'
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws.Range('D13') '0.09371'
Set-TextValue $ws.Range('E13') '-0.28%'
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws.Range('D14') '0.001539'
Set-TextValue $ws.Range('E14') '1.46%'
$ws.Range('B15').Value = 'One'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue $ws.Range('D15') '0.0006032'
Set-TextValue $ws.Range('E15') '-94.09%'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws.Range('D16') '0.005992'
Set-TextValue $ws.Range('E16') '-3.59%'
$ws.Range('B17').Value = 'UpBots'
$ws.Range('C17').Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
Set-TextValue $ws.Range('D17') '0.007489'
Set-TextValue $ws.Range('E17') '5,224.41%'
$ws.Range('B18').Value = 'LEO'
$ws.Range('C18').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range('D18') '3.491'
Set-TextValue $ws.Range('E18') '-0.65%'
Set-TextValue $ws.Range('D19') '2.185'
Set-TextValue $ws.Range('E19') '0.19%'
Set-TextValue $ws.Range('D20') '0.3199'
Set-TextValue $ws.Range('E20') '1.56%'
Set-TextValue $ws.Range('D21') '0.03313'
Set-TextValue $ws.Range('E21') '2.97%'
Set-TextValue $ws.Range('D22') '0.1290'
Set-TextValue $ws.Range('E22') '-1.02%'
Set-TextValue $ws.Range('D23') '3.486'
Set-TextValue $ws.Range('E23') '-2.26%'
Set-TextValue $ws.Range('D24') '0.04155'
Set-TextValue $ws.Range('E24') '1.55%'
Set-TextValue $ws.Range('E25') '0.46%'
Set-TextValue $ws.Range('E26') '1.10%'
Set-TextValue $ws.Range('E27') '11.74%'
Set-TextValue $ws.Range('D28') '0.0001211'
Set-TextValue $ws.Range('E28') '2.63%'
Set-TextValue $ws.Range('D40') '0.03748'
Set-TextValue $ws.Range('E40') '0.84%'
Set-TextValue $ws.Range('D41') '0.003520'
Set-TextValue $ws.Range('E41') '-40.41%'
Set-TextValue $ws.Range('D42') '0.1072'
Set-TextValue $ws.Range('E42') '1.29%'
Set-TextValue $ws.Range('D43') '0.002461'
Set-TextValue $ws.Range('E43') '5.88%'
Set-TextValue $ws.Range('D44') '0.008467'
Set-TextValue $ws.Range('D45') '0.00005284'
Set-TextValue $ws.Range('E45') '3.38%'
Set-TextValue $ws.Range('E46') '0.02%'
Set-TextValue $ws.Range('D47') '0.05701'
Set-TextValue $ws.Range('E47') '-43.56%'
Set-TextValue $ws.Range('D48') '0.002260'
Set-TextValue $ws.Range('E48') '-10.76%'
Set-TextValue $ws.Range('D49') '0.00002100'
Set-TextValue $ws.Range('E49') '0.02%'
Set-TextValue $ws.Range('D50') '0.0002000'
Set-TextValue $ws.Range('E50') '0.02%'
